$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix header text: drop the stray backslash in front of the percent signs ---
$ws.Range("B1").Value = "Difference in mortality rates (%)"
$ws.Range("C1").Value = "95 % CI"

# --- Resize the default/standard column width for the sheet ---
$ws.StandardWidth = 11.53515625

# --- Set explicit widths for the table's columns (A:D) ---
# Column widths are specified in "characters" (same units as the ColumnWidth
# property); the values below are chosen so the resulting stored width is as
# close as possible to the target widths (10.88 / 28.66 / 17.55 / 7.95 chars).
$ws.Columns.Item(1).ColumnWidth = 10
$ws.Columns.Item(2).ColumnWidth = 27.833333333333332
$ws.Columns.Item(3).ColumnWidth = 16.666666666666668
$ws.Columns.Item(4).ColumnWidth = 7.166666666666667

# --- Move the selection / active cell to C1 ---
$ws.Range("C1").Select() | Out-Null
